{"js": "// Change the document's metadata Title from \"SAMPLE AGREEMENTS\" to\n// \"DSRP - Shared Cost Agreement\" (docProps/core.xml <dc:title>).\nconst properties = context.document.properties;\nproperties.load(\"title\");\nawait context.sync();\n\nproperties.title = \"DSRP - Shared Cost Agreement\";\nawait context.sync();\n", "ps1": "# Change the document's metadata Title from \"SAMPLE AGREEMENTS\" to\n# \"DSRP - Shared Cost Agreement\" (docProps/core.xml <dc:title>).\n$d = $word.ActiveDocument\n$d.Title = \"DSRP - Shared Cost Agreement\"\n"}
